$d = $word.ActiveDocument

$replacements = @(
    @{old = "65×47=3055"; new = "33×72=2376"},
    @{old = "78×64=4992"; new = "72×65=4680"},
    @{old = "99×63=6237"; new = "55×86=4730"},
    @{old = "96×51=4896"; new = "45×67=3015"},
    @{old = "37×56=2072"; new = "92×12=1104"},
    @{old = "80×40=3200"; new = "81×31=2511"},
    @{old = "97×97=9409"; new = "79×50=3950"},
    @{old = "48×17=816"; new = "36×93=3348"},
    @{old = "39×72=2808"; new = "41×62=2542"},
    @{old = "15×85=1275"; new = "62×16=992"},
    @{old = "36×11=396"; new = "37×89=3293"},
    @{old = "85×29=2465"; new = "70×12=840"},
    @{old = "42×50=2100"; new = "84×74=6216"},
    @{old = "97×67=6499"; new = "64×92=5888"},
    @{old = "32×40=1280"; new = "18×61=1098"},
    @{old = "88×23=2024"; new = "29×76=2204"},
    @{old = "97×42=4074"; new = "17×44=748"},
    @{old = "15×32=480"; new = "42×15=630"},
    @{old = "43×12=516"; new = "75×53=3975"},
    @{old = "50×70=3500"; new = "73×32=2336"},
    @{old = "50×62=3100"; new = "83×20=1660"},
    @{old = "78×45=3510"; new = "37×18=666"},
    @{old = "81×77=6237"; new = "48×70=3360"},
    @{old = "15×50=750"; new = "58×76=4408"},
    @{old = "42×46=1932"; new = "77×12=924"}
)

foreach ($pair in $replacements) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $true, $false, $false, $false,
                         $true, 1, $false, $pair.new, 2)
}

$d.Save()
